$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,3,4,0),
    @(2,0,4,3),
    @(3,0,6,3),
    @(7,3,4,0),
    @(3,1,4,2),
    @(3,0,2,3),
    @(5,0,2,2),
    @(3,1,3,2),
    @(6,0,5,2),
    @(4,0,3,3),
    @(3,3,2,0),
    @(4,1,6,2),
    @(6,2,6,0),
    @(5,1,4,2),
    @(4,2,3,1),
    @(3,3,3,0),
    @(5,0,6,2),
    @(6,2,6,0),
    @(4,0,3,2),
    @(5,3,5,0),
    @(5,0,5,2),
    @(2,1,6,2),
    @(4,2,4,0),
    @(2,0,4,3),
    @(3,0,3,3),
    @(3,1,4,2),
    @(7,0,7,2),
    @(5,2,5,0),
    @(6,0,6,2),
    @(3,3,3,0),
    @(6,0,6,2),
    @(3,3,4,0),
    @(6,1,6,2),
    @(3,2,4,1),
    @(4,1,5,2),
    @(3,0,2,3),
    @(6,0,3,2),
    @(4,0,4,2),
    @(5,2,4,1),
    @(5,0,5,2),
    @(6,0,5,2),
    @(4,0,4,2),
    @(3,2,3,1),
    @(6,0,5,2),
    @(4,2,4,1),
    @(5,0,7,3)
)

$startRow = 1171
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$ws.Range("E1204").Select()
